$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Change 1: for every data row (2..100), set the "4th & 13/14/15: GFI/FG"
# columns (AL, AN, AO, AQ, AR, AT) from 1.0 down to 0.0 ---
$zeroCols = @("AL","AN","AO","AQ","AR","AT")
for ($r = 2; $r -le 100; $r++) {
    foreach ($col in $zeroCols) {
        $ws.Range("$col$r").Value = 0.0
    }
}

# --- Change 2: for rows 17..23, update the "4th & 7..12" block of values ---
$newVals = @{
    "U"  = 6.0
    "W"  = 2.0
    "X"  = 4.0
    "Y"  = 6.0
    "Z"  = 6.0
    "AA" = 7.0
    "AB" = 8.0
    "AC" = 6.0
    "AD" = 3.0
    "AE" = 6.0
    "AF" = 3.0
    "AG" = 6.0
    "AH" = 2.0
    "AI" = 6.0
}
# these columns also lose their left border as part of the update
$borderClearCols = @("W","Z","AC","AF","AI")

for ($r = 17; $r -le 23; $r++) {
    foreach ($col in $newVals.Keys) {
        $ws.Range("$col$r").Value = $newVals[$col]
    }
    foreach ($col in $borderClearCols) {
        $ws.Range("$col$r").Borders.Item(7).LineStyle = -4142
    }
}
